# Re-upload of the skills-audit sheet: the numbering in column A (rows 12-29)
# was off by one (row 11 and row 12 both showed "5"); renumber rows 12-29
# sequentially as 6-23, and update the sheet's scroll/zoom/selection to match
# where the author left the view (zoomed out further, scrolled down, with
# A29 selected instead of M25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 12; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $cell.Value() + 1
}

$window = $excel.ActiveWindow
$window.ScrollRow = 21
$window.ScrollColumn = 1
$window.Zoom = 55

$ws.Range("A29").Select()
